$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 970; this shifts the existing rows 970-1023
# down to 973-1026 (matching the dimension change A1:R1023 -> A1:R1026).
$ws.Rows("970:972").Insert()

# New date for the three inserted rows: 2022-05-25 (serial 44706)
$newDate = Get-Date -Year 2022 -Month 5 -Day 25 -Hour 0 -Minute 0 -Second 0

# Row 970: Cebollín, Extra
$ws.Cells.Item(970, 1).Value = 9
$ws.Cells.Item(970, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(970, 3).Value = "Metropolitana"
$ws.Cells.Item(970, 4).Value = $newDate
$ws.Cells.Item(970, 5).Value = 13
$ws.Cells.Item(970, 6).Value = 100112037
$ws.Cells.Item(970, 7).Value = "Cebollín"
$ws.Cells.Item(970, 8).Value = "Sin especificar"
$ws.Cells.Item(970, 9).Value = "Extra"
$ws.Cells.Item(970, 10).Value = 180
$ws.Cells.Item(970, 11).Value = 8000
$ws.Cells.Item(970, 12).Value = 8000
$ws.Cells.Item(970, 13).Value = 8000
$ws.Cells.Item(970, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(970, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(970, 16).Value = 222
$ws.Cells.Item(970, 17).Value = 36
$ws.Cells.Item(970, 18).Value = "Hortaliza"

# Row 971: Cebollín, Primera
$ws.Cells.Item(971, 1).Value = 9
$ws.Cells.Item(971, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(971, 3).Value = "Metropolitana"
$ws.Cells.Item(971, 4).Value = $newDate
$ws.Cells.Item(971, 5).Value = 13
$ws.Cells.Item(971, 6).Value = 100112037
$ws.Cells.Item(971, 7).Value = "Cebollín"
$ws.Cells.Item(971, 8).Value = "Sin especificar"
$ws.Cells.Item(971, 9).Value = "Primera"
$ws.Cells.Item(971, 10).Value = 390
$ws.Cells.Item(971, 11).Value = 7000
$ws.Cells.Item(971, 12).Value = 7000
$ws.Cells.Item(971, 13).Value = 7000
$ws.Cells.Item(971, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(971, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(971, 16).Value = 194
$ws.Cells.Item(971, 17).Value = 36
$ws.Cells.Item(971, 18).Value = "Hortaliza"

# Row 972: Cebollín, Segunda
$ws.Cells.Item(972, 1).Value = 9
$ws.Cells.Item(972, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(972, 3).Value = "Metropolitana"
$ws.Cells.Item(972, 4).Value = $newDate
$ws.Cells.Item(972, 5).Value = 13
$ws.Cells.Item(972, 6).Value = 100112037
$ws.Cells.Item(972, 7).Value = "Cebollín"
$ws.Cells.Item(972, 8).Value = "Sin especificar"
$ws.Cells.Item(972, 9).Value = "Segunda"
$ws.Cells.Item(972, 10).Value = 150
$ws.Cells.Item(972, 11).Value = 6000
$ws.Cells.Item(972, 12).Value = 6000
$ws.Cells.Item(972, 13).Value = 6000
$ws.Cells.Item(972, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(972, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(972, 16).Value = 167
$ws.Cells.Item(972, 17).Value = 36
$ws.Cells.Item(972, 18).Value = "Hortaliza"
